$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("G" + $r)
    $val = $cell.Text
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -eq 2 -and $parts[0] -eq "dnasr281@gmail.com") {
            $cell.Value = $parts[1] + ", " + $parts[0]
        }
    }
}
